# Adds columns I ("I0") and J ("IF") to Sheet1, per commit "I0 and IF added".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells: labels + formatting copied from the existing header (H1),
# matching the other header cells (bold, bordered, centered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-52: each entry is (row, I value, J value).
$data = @(
    @(2,8,8),
    @(3,3,4),
    @(4,8,8),
    @(5,6,6),
    @(6,9,9),
    @(7,8,8),
    @(8,8,8),
    @(9,7,8),
    @(10,6,7),
    @(11,7,8),
    @(12,5,5),
    @(13,6,6),
    @(14,10,10),
    @(15,7,7),
    @(16,7,7),
    @(17,7,7),
    @(18,6,7),
    @(19,9,9),
    @(20,8,8),
    @(21,10,11),
    @(22,9,9),
    @(23,6,7),
    @(24,9,9),
    @(25,3,3),
    @(26,4,5),
    @(27,8,8),
    @(28,7,8),
    @(29,6,7),
    @(30,7,7),
    @(31,6,7),
    @(32,5,6),
    @(33,8,8),
    @(34,8,9),
    @(35,4,5),
    @(36,8,8),
    @(37,9,9),
    @(38,6,6),
    @(39,7,8),
    @(40,7,9),
    @(41,9,9),
    @(42,5,7),
    @(43,5,6),
    @(44,4,5),
    @(45,7,8),
    @(46,4,5),
    @(47,4,4),
    @(48,3,3),
    @(49,3,3),
    @(50,6,6),
    @(51,7,7),
    @(52,3,3)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 9).Value = $entry[1]
    $ws.Cells.Item($r, 10).Value = $entry[2]
}

Write-Output "I0/IF columns added for $($data.Count) rows"
